$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Include #0 sheet: update the System URI value ---
$inc0Sheet = $wb.Worksheets.Item("Include #0")
$inc0Sheet.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_G13-OrientationParticuliere/FHIR/TRE-G13-OrientationParticuliere"

# --- Include #1 sheet: update the System URI value ---
$inc1Sheet = $wb.Worksheets.Item("Include #1")
$inc1Sheet.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R01-EnsembleSavoirFaire-CISIS/FHIR/TRE-R01-EnsembleSavoirFaire-CISIS"
